$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, border, centered) from an existing
# header cell onto the three new header cells before setting their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 77
    $ws.Cells.Item($r, 31).Value = 85
    $ws.Cells.Item($r, 32).Value = 0
}
